$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "săpt. 12" (column N) attendance of 2 for the selected students.
$ws.Range("N4").Value = 2
$ws.Range("N7").Value = 2
$ws.Range("N10").Value = 2
$ws.Range("N12").Value = 2
$ws.Range("N13").Value = 2
$ws.Range("N20").Value = 2
$ws.Range("N22").Value = 2

$excel.Calculate()

# Update the active selection in the bottom-right frozen pane.
$ws.Range("N3:N22").Select()
